{"js": "// Highlight the (already-answered) questions in yellow.\n// The author marked each completed question by appending \" (done)\" to the\n// question text; we find those paragraphs and set their highlight color\n// (which also stamps the paragraph mark's run properties, matching how\n// Word applies a \"select whole paragraph + highlight\" action).\nconst body = context.document.body;\nconst doneRanges = body.search(\" (done)\", { matchCase: false, matchWholeWord: false });\ndoneRanges.load(\"items\");\nawait context.sync();\n\nconst targetParagraphs = [];\nfor (let i = 0; i < doneRanges.items.length; i++) {\n  const paras = doneRanges.items[i].paragraphs;\n  paras.load(\"items\");\n  targetParagraphs.push(paras);\n}\nawait context.sync();\n\nfor (let i = 0; i < targetParagraphs.length; i++) {\n  const paras = targetParagraphs[i];\n  for (let j = 0; j < paras.items.length; j++) {\n    paras.items[j].font.highlightColor = \"yellow\";\n  }\n}\nawait context.sync();\n", "ps1": "# Highlight the (already-answered) questions in yellow.\n# The author marked each completed question by appending \" (done)\" to the\n# question text; find those paragraphs and apply a yellow highlight to the\n# whole paragraph (this also stamps the paragraph mark's run properties,\n# matching a \"select whole paragraph + highlight\" action in Word).\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \" (done)\"\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0\n\nwhile ($rng.Find.Execute()) {\n    $para = $rng.Paragraphs(1)\n    $para.Range.Font.HighlightColorIndex = 7\n    $rng.Collapse(0)\n}\n"}
